$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.155.63"
$ws.Range("D2").Style = $s
$s = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("E2").Style = $s
$s = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.34"
$ws.Range("D3").Style = $s
$s = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E3").Style = $s
$s = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("E4").Style = $s
$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.09"
$ws.Range("D5").Style = $s
$s = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("E5").Style = $s
$s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5242"
$ws.Range("D6").Style = $s
$s = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E6").Style = $s
$s = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E7").Style = $s
$s = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2622"
$ws.Range("D8").Style = $s
$s = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06305"
$ws.Range("D9").Style = $s
$s = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("E9").Style = $s
$s = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.60"
$ws.Range("D10").Style = $s
$s = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E10").Style = $s
$s = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07816"
$ws.Range("D11").Style = $s
$s = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("E11").Style = $s
$s = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.499"
$ws.Range("D12").Style = $s
$s = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("E12").Style = $s
$s = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.658.81"
$ws.Range("D13").Style = $s
$s = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("E13").Style = $s
$s = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.884.78"
$ws.Range("D14").Style = $s
$s = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("E14").Style = $s
$s = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5551"
$ws.Range("D15").Style = $s
$s = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("E15").Style = $s
$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8025"
$ws.Range("D16").Style = $s
$s = $ws.Range("E16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("E16").Style = $s
$s = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.09"
$ws.Range("D17").Style = $s
$s = $ws.Range("E17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("E17").Style = $s
$s = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.158.15"
$ws.Range("D18").Style = $s
$s = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("E18").Style = $s
$s = $ws.Range("E19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("E19").Style = $s
$s = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.645"
$ws.Range("D20").Style = $s
$s = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.43"
$ws.Range("D21").Style = $s
$s = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("E21").Style = $s
$s = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("E22").Style = $s
$s = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.962"
$ws.Range("D23").Style = $s
$s = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("E23").Style = $s
$s = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("E24").Style = $s
$s = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.56"
$ws.Range("D25").Style = $s
$s = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("E25").Style = $s
$s = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("E26").Style = $s
$s = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.153"
$ws.Range("D27").Style = $s
$s = $ws.Range("E27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("E27").Style = $s
$s = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.96"
$ws.Range("D28").Style = $s
$s = $ws.Range("E28").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("E28").Style = $s
$s = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.493"
$ws.Range("D29").Style = $s
$s = $ws.Range("E29").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E29").Style = $s
$s = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05710"
$ws.Range("D30").Style = $s
$s = $ws.Range("E30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.09%  "
$ws.Range("E30").Style = $s
$s = $ws.Range("E31").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("E31").Style = $s
$s = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.486"
$ws.Range("D32").Style = $s
$s = $ws.Range("E32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("E32").Style = $s
$s = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.356"
$ws.Range("D33").Style = $s
$s = $ws.Range("E34").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("E34").Style = $s
$s = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9538"
$ws.Range("D35").Style = $s
$s = $ws.Range("E35").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("E35").Style = $s
$s = $ws.Range("E36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("E36").Style = $s
$s = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.422"
$ws.Range("D37").Style = $s
$s = $ws.Range("E37").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E37").Style = $s
$s = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5697"
$ws.Range("D38").Style = $s
$s = $ws.Range("E38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("E38").Style = $s
$s = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01595"
$ws.Range("D39").Style = $s
$s = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("E39").Style = $s
$s = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.943"
$ws.Range("D40").Style = $s
$s = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("E40").Style = $s
$s = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.063.98"
$ws.Range("D41").Style = $s
$s = $ws.Range("E41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("E41").Style = $s
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("D42").Style = $s
$s = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("E42").Style = $s
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$s = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8439"
$ws.Range("D43").Style = $s
$s = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("E43").Style = $s
$s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.51"
$ws.Range("D44").Style = $s
$s = $ws.Range("E44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("E44").Style = $s
$s = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.795.84"
$ws.Range("D45").Style = $s
$s = $ws.Range("E45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("E45").Style = $s
$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.87"
$ws.Range("D46").Style = $s
$s = $ws.Range("E46").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("E46").Style = $s
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$s = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05453"
$ws.Range("D47").Style = $s
$s = $ws.Range("E47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.64%  "
$ws.Range("E47").Style = $s
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$s = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").Style = $s
$s = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("E48").Style = $s
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$s = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4403"
$ws.Range("D49").Style = $s
$s = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E49").Style = $s
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$s = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.973"
$ws.Range("D50").Style = $s
$s = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("E50").Style = $s
$s = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.42%  "
$ws.Range("E51").Style = $s
